# Applies the "cbc -> cplex" solver rename + scenario-label split, plus the
# small Setup-sheet value tweaks described by the commit's XML diff.

$wb = $excel.ActiveWorkbook

$wsUncertainty = $wb.Worksheets.Item("Uncertainty_Table")
$wsSetup       = $wb.Worksheets.Item("Setup")

# --- Setup sheet: solver cbc -> cplex (must precede the Scenario1 split --
#     below so the shared-string table grows in the same order Excel used)
$wsSetup.Range("H2").Value = "cplex"
$wsSetup.Range("A2").Value = 2
$wsSetup.Range("D2").Value = 2

# --- Uncertainty_Table: "BAU ; Scenario1" -> "Scenario1" for D2:D15 -------
$rng = $wsUncertainty.Range("D2:D15")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq "BAU ; Scenario1") {
        $cell.Value = "Scenario1"
    }
}

# --- Make "Setup" the active/selected sheet & tab, matching the new view -
$wsUncertainty.Range("C13").Select()
$wsUncertainty.Application.ActiveWindow.ScrollRow = 11

$wsSetup.Activate()
$wsSetup.Range("F3").Select()
